$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.038140866589673
$ws.Range("D2").Value = 1.035478105637088
$ws.Range("E2").Value = 1.036713933457367
$ws.Range("F2").Value = 1.036855476723082
$ws.Range("I2").Value = 1.038314527069739
$ws.Range("J2").Value = 1.043239918403781
$ws.Range("K2").Value = 1.038274595874669
$ws.Range("L2").Value = 1.039506886416251
$ws.Range("M2").Value = 1.039648025130748
$ws.Range("N2").Value = 1.044721438749274

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.03986110267309
$ws.Range("D3").Value = 1.036208875121174
$ws.Range("E3").Value = 1.038207001999305
$ws.Range("F3").Value = 1.039213542259758
$ws.Range("I3").Value = 1.038769739696946
$ws.Range("J3").Value = 1.044600960472505
$ws.Range("K3").Value = 1.03881497757039
$ws.Range("L3").Value = 1.040807802397683
$ws.Range("M3").Value = 1.041811680137632
$ws.Range("N3").Value = 1.046084413653781

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.040970149254108
$ws.Range("D4").Value = 1.036680132471149
$ws.Range("E4").Value = 1.039169558261267
$ws.Range("F4").Value = 1.040734650820889
$ws.Range("I4").Value = 1.03906133762111
$ws.Range("J4").Value = 1.045477326239966
$ws.Range("K4").Value = 1.03916237839967
$ws.Range("L4").Value = 1.04164553380857
$ws.Range("M4").Value = 1.043206700744708
$ws.Range("N4").Value = 1.046962023961152

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.041435439614219
$ws.Range("D5").Value = 1.036877870857753
$ws.Range("E5").Value = 1.039573380853299
$ws.Range("F5").Value = 1.041373027166922
$ws.Range("I5").Value = 1.039183223543237
$ws.Range("J5").Value = 1.045844732373075
$ws.Range("K5").Value = 1.039307888972226
$ws.Range("L5").Value = 1.041996761908095
$ws.Range("M5").Value = 1.043791997588426
$ws.Range("N5").Value = 1.047329951853113

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.041513508599431
$ws.Range("D6").Value = 1.036911049959719
$ws.Range("E6").Value = 1.039641135876314
$ws.Range("F6").Value = 1.04148014990875
$ws.Range("I6").Value = 1.039203647730813
$ws.Range("J6").Value = 1.045906362244042
$ws.Range("K6").Value = 1.039332289502843
$ws.Range("L6").Value = 1.042055679128973
$ws.Range("M6").Value = 1.043890203750106
$ws.Range("N6").Value = 1.047391669245565

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.04097637020426
$ws.Range("D7").Value = 1.03668277614072
$ws.Range("E7").Value = 1.039174957421977
$ws.Range("F7").Value = 1.040743185107367
$ws.Range("I7").Value = 1.039062969016872
$ws.Range("J7").Value = 1.0454822395171
$ws.Range("K7").Value = 1.039164324821388
$ws.Range("L7").Value = 1.041650230664598
$ws.Range("M7").Value = 1.043214526066867
$ws.Range("N7").Value = 1.046966944215703

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.038723081090047
$ws.Range("D8").Value = 1.035725406611352
$ws.Range("E8").Value = 1.037219269517685
$ws.Range("F8").Value = 1.037653391989859
$ws.Range("I8").Value = 1.038468983576263
$ws.Range("J8").Value = 1.043700793035233
$ws.Range("K8").Value = 1.038457691757766
$ws.Range("L8").Value = 1.039947384334199
$ws.Range("M8").Value = 1.040380297333382
$ws.Range("N8").Value = 1.045182967875497

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.034720514914967
$ws.Range("D9").Value = 1.034025936397334
$ws.Range("E9").Value = 1.033745118796011
$ws.Range("F9").Value = 1.032171218146889
$ws.Range("I9").Value = 1.037399418563429
$ws.Range("J9").Value = 1.040527848416889
$ws.Range("K9").Value = 1.037194967300348
$ws.Range("L9").Value = 1.036915069874215
$ws.Range("M9").Value = 1.035346336699589
$ws.Range("N9").Value = 1.042005517311853

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.032029395828266
$ws.Range("D10").Value = 1.032884273184938
$ws.Range("E10").Value = 1.031409189868671
$ws.Range("F10").Value = 1.028489162904156
$ws.Range("I10").Value = 1.036670627673494
$ws.Range("J10").Value = 1.03838879719825
$ws.Range("K10").Value = 1.036341040944949
$ws.Range("L10").Value = 1.034871271455367
$ws.Range("M10").Value = 1.031961813214676
$ws.Range("N10").Value = 1.039863428395131

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.030858439388013
$ws.Range("D11").Value = 1.032387791739671
$ws.Range("E11").Value = 1.030392775157118
$ws.Range("F11").Value = 1.026887870972091
$ws.Range("I11").Value = 1.03635123678712
$ws.Range("J11").Value = 1.037456699806287
$ws.Range("K11").Value = 1.03596833821108
$ws.Range("L11").Value = 1.033980790861792
$ws.Range("M11").Value = 1.030489086703723
$ws.Range("N11").Value = 1.038930007317958

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.030422615904799
$ws.Range("D12").Value = 1.032203049763995
$ws.Range("E12").Value = 1.030014471116093
$ws.Range("F12").Value = 1.026291998319505
$ws.Range("I12").Value = 1.036232019821076
$ws.Range("J12").Value = 1.03710957479424
$ws.Range("K12").Value = 1.035829450850619
$ws.Range("L12").Value = 1.033649181387268
$ws.Range("M12").Value = 1.029940931651715
$ws.Range("N12").Value = 1.038582389348579

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.030516141649932
$ws.Range("D13").Value = 1.032242692426004
$ws.Range("E13").Value = 1.030095653437733
$ws.Range("F13").Value = 1.026419864603886
$ws.Range("I13").Value = 1.036257618670562
$ws.Range("J13").Value = 1.037184075454158
$ws.Range("K13").Value = 1.035859263078752
$ws.Range("L13").Value = 1.033720351300946
$ws.Range("M13").Value = 1.030058563990414
$ws.Range("N13").Value = 1.038656995807979

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.030822432137432
$ws.Range("D14").Value = 1.032372527609771
$ws.Range("E14").Value = 1.030361520095057
$ws.Range("F14").Value = 1.026838638244761
$ws.Range("I14").Value = 1.036341394168713
$ws.Range("J14").Value = 1.037428024872973
$ws.Range("K14").Value = 1.03595686693377
$ws.Range("L14").Value = 1.033953397272713
$ws.Range("M14").Value = 1.030443799072472
$ws.Range("N14").Value = 1.038901291662951

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.031011030726859
$ws.Range("D15").Value = 1.032452479881623
$ws.Range("E15").Value = 1.030525227697402
$ws.Range("F15").Value = 1.027096514008889
$ws.Range("I15").Value = 1.036392933819024
$ws.Range("J15").Value = 1.037578209961749
$ws.Range("K15").Value = 1.036016944209304
$ws.Range("L15").Value = 1.034096872009328
$ws.Range("M15").Value = 1.030681005766087
$ws.Range("N15").Value = 1.039051690031778

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.032106986539879
$ws.Range("D16").Value = 1.032917177550639
$ws.Range("E16").Value = 1.0314765401526
$ws.Range("F16").Value = 1.02859528587273
$ws.Range("I16").Value = 1.036691743515063
$ws.Range("J16").Value = 1.038450531969884
$ws.Range("K16").Value = 1.036365713379944
$ws.Range("L16").Value = 1.034930252204207
$ws.Range("M16").Value = 1.03205939832505
$ws.Range("N16").Value = 1.039925250837222

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.032792911931977
$ws.Range("D17").Value = 1.033208094594591
$ws.Range("E17").Value = 1.032071936255768
$ws.Range("F17").Value = 1.029533540948186
$ws.Range("I17").Value = 1.036878151335535
$ws.Range("J17").Value = 1.038996129817897
$ws.Range("K17").Value = 1.036583693717741
$ws.Range("L17").Value = 1.035451523351352
$ws.Range("M17").Value = 1.032922072988084
$ws.Range("N17").Value = 1.04047162349676

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.033192453822811
$ws.Range("D18").Value = 1.033377576127311
$ws.Range("E18").Value = 1.032418745248657
$ws.Range("F18").Value = 1.030080142549698
$ws.Range("I18").Value = 1.036986511812097
$ws.Range("J18").Value = 1.039313802414476
$ws.Range("K18").Value = 1.036710554045313
$ws.Range("L18").Value = 1.035755042366173
$ws.Range("M18").Value = 1.033424563776726
$ws.Range("N18").Value = 1.040789747224864

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.033328595263534
$ws.Range("D19").Value = 1.033435330283418
$ws.Range("E19").Value = 1.032536918099819
$ws.Range("F19").Value = 1.03026640766266
$ws.Range("I19").Value = 1.037023397719867
$ws.Range("J19").Value = 1.039422025260514
$ws.Range("K19").Value = 1.036753762168058
$ws.Range("L19").Value = 1.035858445233196
$ws.Range("M19").Value = 1.033595783738619
$ws.Range("N19").Value = 1.040898123759756

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.032719375359393
$ws.Range("D20").Value = 1.033176903244244
$ws.Range("E20").Value = 1.032008105212961
$ws.Range("F20").Value = 1.029432944301221
$ws.Range("I20").Value = 1.036858189661046
$ws.Range("J20").Value = 1.038937650958133
$ws.Range("K20").Value = 1.03656033590194
$ws.Range("L20").Value = 1.035395650806596
$ws.Range("M20").Value = 1.032829588076285
$ws.Range("N20").Value = 1.040413061590308

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.030732261659236
$ws.Range("D21").Value = 1.032334303451139
$ws.Range("E21").Value = 1.030283250197005
$ws.Range("F21").Value = 1.026715350018734
$ws.Range("I21").Value = 1.036316740467892
$ws.Range("J21").Value = 1.037356212926229
$ws.Range("K21").Value = 1.035928137475445
$ws.Range("L21").Value = 1.033884794529162
$ws.Range("M21").Value = 1.030330388122724
$ws.Range("N21").Value = 1.038829377735005

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.029477791621856
$ws.Range("D22").Value = 1.031802634456791
$ws.Range("E22").Value = 1.02919434453885
$ws.Range("F22").Value = 1.025000413350031
$ws.Range("I22").Value = 1.035972944946994
$ws.Range("J22").Value = 1.036356668112768
$ws.Range("K22").Value = 1.035528048467993
$ws.Range("L22").Value = 1.032929958931746
$ws.Range("M22").Value = 1.028752550183945
$ws.Range("N22").Value = 1.03782841345326

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.030143300819242
$ws.Range("D23").Value = 1.032084663754038
$ws.Range("E23").Value = 1.029772019789404
$ws.Range("F23").Value = 1.025910141993957
$ws.Range("I23").Value = 1.036155518819281
$ws.Range("J23").Value = 1.036887048424612
$ws.Range("K23").Value = 1.035740391847776
$ws.Range("L23").Value = 1.033436606143556
$ws.Range("M23").Value = 1.029589619540451
$ws.Range("N23").Value = 1.038359546965982

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.03275260506006
$ws.Range("D24").Value = 1.033190997907726
$ws.Range("E24").Value = 1.032036949182338
$ws.Range("F24").Value = 1.029478401654522
$ws.Range("I24").Value = 1.036867210619613
$ws.Range("J24").Value = 1.038964076784158
$ws.Range("K24").Value = 1.036570891170693
$ws.Range("L24").Value = 1.035420898840582
$ws.Range("M24").Value = 1.032871380166205
$ws.Range("N24").Value = 1.040439524944037

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.035759190731099
$ws.Range("D25").Value = 1.034466798006913
$ws.Range("E25").Value = 1.034646689471459
$ws.Range("F25").Value = 1.033593148824779
$ws.Range("I25").Value = 1.037678675637715
$ws.Range("J25").Value = 1.041352241611754
$ws.Range("K25").Value = 1.037523523068452
$ws.Range("L25").Value = 1.037702847698863
$ws.Range("M25").Value = 1.042831081239615
$ws.Range("N25").Value = 1.042831081239615
